$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 555

$ws.Range("E3").Value = 425

$ws.Range("E4").Value = 422

$ws.Range("A5").Value = "cand_22_cv_02094"
$ws.Range("D5").Value = 75198
$ws.Range("E5").Value = 453

$ws.Range("A6").Value = "cand_23_cv_02560"
$ws.Range("D6").Value = 65569
$ws.Range("E6").Value = 403

$ws.Range("A7").Value = "cand_23_cv_03518"
$ws.Range("D7").Value = 30953
$ws.Range("E7").Value = 688

$ws.Range("A8").Value = "cand_24_cv_03170"
$ws.Range("D8").Value = 25000
$ws.Range("E8").Value = 586

$ws.Range("A9").Value = "cand_24_cv_04196"
$ws.Range("D9").Value = 11840
$ws.Range("E9").Value = 435

$ws.Range("A10").Value = "cand_3_22-cv-00956"
$ws.Range("D10").Value = 20290
$ws.Range("E10").Value = 378

$ws.Range("A11").Value = "casd_3_23-cv-01216"
$ws.Range("D11").Value = 33148
$ws.Range("E11").Value = 464

$ws.Range("A12").Value = "ctd-3-23-cv-01035"
$ws.Range("D12").Value = 62281
$ws.Range("E12").Value = 324

$ws.Range("A13").Value = "dcd-1_23-cv-02055"
$ws.Range("D13").Value = 36392
$ws.Range("E13").Value = 418

$ws.Range("A14").Value = "dde_ 23_cv_1466"
$ws.Range("D14").Value = 33685
$ws.Range("E14").Value = 462

$ws.Range("A15").Value = "dde_21_cv_55"
$ws.Range("D15").Value = 44368
$ws.Range("E15").Value = 369

$ws.Range("A16").Value = "flsd-1_23-cv-23139"
$ws.Range("D16").Value = 15671
$ws.Range("E16").Value = 885

$ws.Range("A17").Value = "ilnd-1-21-cv-04349"
$ws.Range("D17").Value = 30291
$ws.Range("E17").Value = 733

$ws.Range("A18").Value = "mad-1-21-cv-10933"
$ws.Range("D18").Value = 16572
$ws.Range("E18").Value = 421

$ws.Range("A19").Value = "mied-4-23-cv-13132"
$ws.Range("D19").Value = 64016
$ws.Range("E19").Value = 552

$ws.Range("A20").Value = "nysd_20_cv_04494"
$ws.Range("D20").Value = 50338
$ws.Range("E20").Value = 535

$ws.Range("A21").Value = "nysd_22-cv-07111"
$ws.Range("D21").Value = 27127
$ws.Range("E21").Value = 569

$ws.Range("A22").Value = "nysd_22_cv_10292"
$ws.Range("D22").Value = 23587
$ws.Range("E22").Value = 411

$ws.Range("A23").Value = "nysd_23_cv_9476"
$ws.Range("D23").Value = 16095
$ws.Range("E23").Value = 415

$ws.Range("A24").Value = "nysd_24_cv_310"
$ws.Range("D24").Value = 43285
$ws.Range("E24").Value = 546

$ws.Range("A25").Value = "txnd-4_24-cv-00673"
$ws.Range("D25").Value = 54955
$ws.Range("E25").Value = 428

$ws.Range("A26").Value = "txsd-4-21-cv-02473"
$ws.Range("D26").Value = 66230
$ws.Range("E26").Value = 444
